$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ideas")

$ws.Range("C14").Value = "Social Network"
$ws.Range("D14").Value = "Social Network"

$ws.Range("D14").Select()
